$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update rows 16-30: flip sign of values in columns A, C, D, E, F (and adjust C from 3 to 2)
for ($r = 16; $r -le 30; $r++) {
    $ws.Cells.Item($r, 1).Value = 1000
    $ws.Cells.Item($r, 3).Value = 2
    $ws.Cells.Item($r, 4).Value = 50
    $ws.Cells.Item($r, 5).Value = 500
    $ws.Cells.Item($r, 6).Value = 150
}

# Delete rows 31-60 entirely (shrinks used range to A1:G30)
$ws.Rows("31:60").Delete()

# Update the selection to match the new state
$ws.Range("A31:G31").Select()
